# Re-order the rows of the "classFields" sheet so that, within each class,
# fields appear in the order produced by the regenerated structure report:
#   - pl.piomin.stock.domain.Product: rows rotate from
#       [availableItems, name, id, reservedItems]
#     to
#       [reservedItems, id, name, availableItems]
#   - pl.piomin.stock.StockComponentTests: rows rotate from
#       [template, product, repository, factory, LOG, kafka]
#     to
#       [product, LOG, repository, kafka, template, factory]
#   - pl.piomin.stock.service.OrderManageService: rows rotate from
#       [repository, SOURCE, LOG, template]
#     to
#       [repository, SOURCE, template, LOG]
# The set of (class, field, modifier, type) rows is unchanged overall; only
# the row order (and hence which row each field lands on) changes. Only the
# cells whose resolved text actually changes are touched below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# --- pl.piomin.stock.domain.Product (rows 2-5) ---
$ws.Cells.Item(2, 2).Value = "reservedItems"

$ws.Cells.Item(3, 2).Value = "id"
$ws.Cells.Item(3, 4).Value = "java.lang.Long"

$ws.Cells.Item(4, 2).Value = "name"
$ws.Cells.Item(4, 4).Value = "java.lang.String"

$ws.Cells.Item(5, 2).Value = "availableItems"

# --- pl.piomin.stock.StockComponentTests (rows 6-11) ---
$ws.Cells.Item(6, 2).Value = "product"
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = "pl.piomin.stock.domain.Product"

$ws.Cells.Item(7, 2).Value = "LOG"
$ws.Cells.Item(7, 3).Value = "private"
$ws.Cells.Item(7, 4).Value = "org.slf4j.Logger"

# row 8 ("repository", package-private, ProductRepository) is unchanged

$ws.Cells.Item(9, 2).Value = "kafka"
$ws.Cells.Item(9, 4).Value = "org.springframework.kafka.test.EmbeddedKafkaBroker"

$ws.Cells.Item(10, 2).Value = "template"
$ws.Cells.Item(10, 4).Value = "org.springframework.kafka.core.KafkaTemplate"

$ws.Cells.Item(11, 2).Value = "factory"
$ws.Cells.Item(11, 4).Value = "org.springframework.kafka.core.ConsumerFactory"

# --- pl.piomin.stock.service.OrderManageService (rows 12-15) ---
# rows 12 ("repository") and 13 ("SOURCE") are unchanged

$ws.Cells.Item(14, 2).Value = "template"
$ws.Cells.Item(14, 4).Value = "org.springframework.kafka.core.KafkaTemplate"

$ws.Cells.Item(15, 2).Value = "LOG"
$ws.Cells.Item(15, 4).Value = "org.slf4j.Logger"

# --- pl.piomin.stock.StockApp (rows 16-18) are unchanged ---
